# Update countries & provincias Spain
# - Refresh the "last updated" timestamp (19:22 -> 19:52).
# - Ecuador's case counts overtook Arabia Saudita and Polonia, so it now
#   ranks just below Chile: row 29 becomes Ecuador (new data), and the
#   previous Arabia Saudita / Polonia rows shift down into rows 30 / 31
#   (each keeping its own data, just one rank lower).
# - A handful of other countries (rows 4, 7, 53, 108, 110, 162) got refreshed
#   case/recovery/death counts from the same data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 19:52"

# --- Estados Unidos (row 4) refreshed figures ------------------------------
$ws.Range("B4").Value = 755162
$ws.Range("C4").Value = 16370
$ws.Range("E4").Value = 646231
$ws.Range("G4").Value = 1095
$ws.Range("H4").Value = 40109

# --- Reino Unido (row 7) refreshed figures ---------------------------------
$ws.Range("D7").Value = 36578
$ws.Range("E7").Value = 96282

# --- Ranking shuffle around Chile: Ecuador, Arabia Saudita, Polonia --------
# Row 29: now Ecuador, with fresh data.
$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 9468
$ws.Range("C29").Value = 446
$ws.Range("D29").Value = 1061
$ws.Range("E29").Value = 7933
$ws.Range("F29").Value = 124
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = 474

# Row 30: now Arabia Saudita (previous row-29 data, one rank down).
$ws.Range("A30").Value = "Arabia Saudita"
$ws.Range("B30").Value = 9362
$ws.Range("C30").Value = 1088
$ws.Range("D30").Value = 1398
$ws.Range("E30").Value = 7867
$ws.Range("F30").Value = 97
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 97

# Row 31: now Polonia (previous row-30 data, one rank down).
$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 9287
$ws.Range("C31").Value = 545
$ws.Range("D31").Value = 1040
$ws.Range("E31").Value = 7887
$ws.Range("F31").Value = 160
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 360

# --- Sudafrica (row 53) refreshed figures ----------------------------------
$ws.Range("B53").Value = 3158
$ws.Range("C53").Value = 124
$ws.Range("E53").Value = 2201
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 54

# --- Jordania (row 108) refreshed figures ----------------------------------
$ws.Range("B108").Value = 417
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 276
$ws.Range("E108").Value = 134

# --- Georgia (row 110) refreshed figures -----------------------------------
$ws.Range("D110").Value = 93
$ws.Range("E110").Value = 297

# --- Macao (row 162) refreshed figures --------------------------------------
$ws.Range("D162").Value = 17
$ws.Range("E162").Value = 28
